# "temp table berubah bknya" - the new master-data row ("master Hotel" /
# "MST020") was actually typed on the "Master" sheet, not "Temp". Add the
# row there, make Master the active/selected sheet (which also clears the
# stale tabSelected flag left on "Temp"), and leave the selection one row
# below the freshly entered data - exactly like Excel does after you type
# a value into a cell and press Enter.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master")

$ws.Activate()

$ws.Range("A8").Value = "master Hotel"
$ws.Range("B8").Value = "MST020"

$ws.Range("A9").Select()
